# Update the "dSF" column (F) values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -2
    4  = -4
    6  = 4
    8  = 11
    9  = 4
    10 = -4
    11 = 2
    12 = -3
    13 = 5
    14 = -2
    15 = -2
    16 = 2
    17 = 6
    18 = -4
    19 = 5
    20 = 1
    21 = -4
    22 = 1
    23 = -7
    25 = 3
    26 = 1
    28 = -2
    29 = 1
    31 = 2
    32 = -3
    33 = -1
    34 = -1
    35 = 3
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
